$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the style/format from the
# neighboring header cell (G1) and then overwriting its value/text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column (H2:H8) with the corresponding values.
$saveValues = @(0, 0, 1, 1, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
